# Updated Global_M2 for easier usage.
# Applies corrected values to rows 253-255 and appends new rows 256-258.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct previously-existing rows (revised source values) ---

# Row 253
$ws.Range("C253").Value = 5817095330000
$ws.Range("D253").Value = 5817095330000
$ws.Range("E253").Value = 5817095330000
$ws.Range("F253").Value = 5817095330000

# Row 254
$ws.Range("C254").Value = 5950864520000
$ws.Range("D254").Value = 5950864520000
$ws.Range("E254").Value = 5950864520000
$ws.Range("F254").Value = 5950864520000

# Row 255
$ws.Range("C255").Value = 6045092150000
$ws.Range("D255").Value = 6045092150000
$ws.Range("E255").Value = 6045092150000
$ws.Range("F255").Value = 6045092150000

# --- Append new rows 256-258 (same layout/style as the existing data rows) ---

$newRows = @(
    @{ Row = 256; Date = 44986.45833333334; Value = 6077620130000 },
    @{ Row = 257; Date = 45017.45833333334; Value = 6141246740000 },
    @{ Row = 258; Date = 45047.41666666666; Value = 6224272840000 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy formatting (style) from the row above so the new row matches
    # the existing data rows (bordered, centered, date-formatted column A).
    $ws.Range("A" + ($row - 1) + ":G" + ($row - 1)).Copy()
    $ws.Range("A" + $row + ":G" + $row).PasteSpecial(-4122) # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "ECONOMICS:CZM2"
    $ws.Cells.Item($row, 3).Value = $r.Value
    $ws.Cells.Item($row, 4).Value = $r.Value
    $ws.Cells.Item($row, 5).Value = $r.Value
    $ws.Cells.Item($row, 6).Value = $r.Value
    $ws.Cells.Item($row, 7).Value = 0
}

$excel.CutCopyMode = $false
